# Updates the division problems in the single table of the worksheet.
# Each data row of the table (rows 1, 5, 9, 13, 17 -- interleaved with
# blank spacer rows) holds 5 cells of the form "NN÷N=". We replace the
# old expression in each cell with the new one. The replacement is
# scoped to the individual cell's Range (with Wrap:=wdFindStop (0) and
# Replace:=wdReplaceOne (1)) so that duplicate text elsewhere in the
# document (e.g. "81÷9=" appears twice) is not clobbered.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wdFindStop = 0
$wdReplaceOne = 1

$changes = @(
    @{ Row = 1;  Col = 1; Old = "81÷9="; New = "82÷7=" },
    @{ Row = 1;  Col = 2; Old = "70÷6="; New = "91÷8=" },
    @{ Row = 1;  Col = 3; Old = "33÷9="; New = "37÷4=" },
    @{ Row = 1;  Col = 4; Old = "16÷8="; New = "85÷5=" },
    @{ Row = 1;  Col = 5; Old = "47÷9="; New = "52÷9=" },

    @{ Row = 5;  Col = 1; Old = "46÷4="; New = "45÷8=" },
    @{ Row = 5;  Col = 2; Old = "37÷8="; New = "75÷6=" },
    @{ Row = 5;  Col = 3; Old = "37÷5="; New = "38÷4=" },
    @{ Row = 5;  Col = 4; Old = "18÷8="; New = "15÷2=" },
    @{ Row = 5;  Col = 5; Old = "85÷7="; New = "88÷9=" },

    @{ Row = 9;  Col = 1; Old = "90÷9="; New = "69÷3=" },
    @{ Row = 9;  Col = 2; Old = "83÷8="; New = "77÷3=" },
    @{ Row = 9;  Col = 3; Old = "18÷9="; New = "13÷7=" },
    @{ Row = 9;  Col = 4; Old = "15÷8="; New = "66÷6=" },
    @{ Row = 9;  Col = 5; Old = "98÷4="; New = "73÷3=" },

    @{ Row = 13; Col = 1; Old = "91÷3="; New = "19÷4=" },
    @{ Row = 13; Col = 2; Old = "30÷6="; New = "57÷9=" },
    @{ Row = 13; Col = 3; Old = "81÷9="; New = "40÷9=" },
    @{ Row = 13; Col = 4; Old = "63÷8="; New = "44÷6=" },
    @{ Row = 13; Col = 5; Old = "21÷6="; New = "88÷9=" },

    @{ Row = 17; Col = 1; Old = "41÷8="; New = "67÷7=" },
    @{ Row = 17; Col = 2; Old = "43÷6="; New = "37÷3=" },
    @{ Row = 17; Col = 3; Old = "46÷8="; New = "50÷5=" },
    @{ Row = 17; Col = 4; Old = "86÷4="; New = "43÷8=" },
    @{ Row = 17; Col = 5; Old = "81÷5="; New = "99÷8=" }
)

foreach ($change in $changes) {
    $cell = $t.Cell($change.Row, $change.Col)
    $rng = $cell.Range
    $rng.Find.Execute($change.Old, $true, $false, $false, $false, $false, `
                       $true, $wdFindStop, $false, $change.New, $wdReplaceOne)
}
